{"js": "// Add \"other issues and command\" right before the existing _GoBack\n// bookmark (same paragraph), then add a new paragraph after it containing\n// the sqlshack URL as plain text.\nconst doc = context.document;\n\n// Locate the paragraph holding the \"_GoBack\" bookmark and insert the new\n// run of text immediately before the bookmark start.\nconst bookmarkRange = doc.getBookmarkRange(\"_GoBack\");\nbookmarkRange.insertText(\"other issues and command\", \"Before\");\n\n// Append a brand-new paragraph at the very end of the body with the URL.\ncontext.document.body.insertParagraph(\n  \"https://www.sqlshack.com/azure-kubernetes-service-aks-managing-sql-server-database-files/\",\n  \"End\"\n);\n\nawait context.sync();\n", "ps1": "# Add \"other issues and command\" right before the existing _GoBack bookmark\n# (staying inside that same paragraph), then add a brand-new paragraph at\n# the end of the document containing the sqlshack URL as plain text.\n\n$d = $word.ActiveDocument\n\n# Paragraph that holds the \"_GoBack\" bookmark: insert the new run of text\n# immediately before the bookmark.\n$bm = $d.Bookmarks(\"_GoBack\")\n$bm.Range.InsertBefore(\"other issues and command\")\n\n# Append a new, final paragraph with the URL text.\n$lastPara = $d.Paragraphs.Last\n$lastPara.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.InsertBefore(\"https://www.sqlshack.com/azure-kubernetes-service-aks-managing-sql-server-database-files/\")\n"}
